$wb = $excel.ActiveWorkbook

# Update info_solution!A2 (comp_time)
$wsInfo = $wb.Worksheets.Item("info_solution")
$wsInfo.Range("A2").Value = 3.638042449951172

# Update design_users table values
$wsDesign = $wb.Worksheets.Item("design_users")
$wsDesign.Range("B2").Value = 30523.270014740006
$wsDesign.Range("C2").Value = 232.10403242918812
$wsDesign.Range("B3").Value = 12160.103553284996
$wsDesign.Range("C3").Value = 92.467454113238
$wsDesign.Range("B4").Value = 24500.351344201987
$wsDesign.Range("C4").Value = 186.30475503361959
